$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15727
$ws.Range("D2").Value = 6587
$ws.Range("E2").Value = 24969039
$ws.Range("C4").Value = 30720
$ws.Range("E4").Value = 47086571
$ws.Range("C7").Value = 72891
$ws.Range("D7").Value = 31023
$ws.Range("E7").Value = 111313801
$ws.Range("C8").Value = 119255
$ws.Range("D8").Value = 38083
$ws.Range("E8").Value = 251695615
$ws.Range("C9").Value = 44146
$ws.Range("D9").Value = 9557
$ws.Range("E9").Value = 85552253
$ws.Range("D10").Value = 27728
$ws.Range("C11").Value = 15168
$ws.Range("E11").Value = 26099684
$ws.Range("C14").Value = 72715
$ws.Range("D14").Value = 23569
$ws.Range("E14").Value = 136881268
$ws.Range("C15").Value = 34927
$ws.Range("D15").Value = 10657
$ws.Range("E15").Value = 67247816
$ws.Range("C17").Value = 72189
$ws.Range("D17").Value = 18315
$ws.Range("E17").Value = 112490146
$ws.Range("C19").Value = 58998
$ws.Range("E19").Value = 114283809
$ws.Range("C20").Value = 80430
$ws.Range("E20").Value = 138918240
$ws.Range("C28").Value = 4709
$ws.Range("E28").Value = 8831038
$ws.Range("C29").Value = 37013
$ws.Range("E29").Value = 118452912
$ws.Range("C33").Value = 13782
$ws.Range("E33").Value = 26325441
$ws.Range("E36").Value = 15575061
$ws.Range("C42").Value = 11250
$ws.Range("E42").Value = 16221836
$ws.Range("C46").Value = 38111
$ws.Range("E46").Value = 81729819
$ws.Range("C47").Value = 6475
$ws.Range("D47").Value = 1798
$ws.Range("E47").Value = 12028650
$ws.Range("C52").Value = 20598
$ws.Range("E52").Value = 40226132
$ws.Range("C53").Value = 8307
$ws.Range("E53").Value = 15115979
$ws.Range("C54").Value = 11390
$ws.Range("E54").Value = 18069636
$ws.Range("C56").Value = 17692
$ws.Range("E56").Value = 32867524
$ws.Range("C69").Value = 12949
$ws.Range("E69").Value = 24094501
$ws.Range("C87").Value = 1738
$ws.Range("E87").Value = 2611989
$ws.Range("C100").Value = 5995
$ws.Range("E100").Value = 10224165
$ws.Range("C104").Value = 13531
$ws.Range("E104").Value = 25029777
$ws.Range("C106").Value = 17058
$ws.Range("E106").Value = 26871140
$ws.Range("C107").Value = 23958
$ws.Range("E107").Value = 32585291
$ws.Range("C108").Value = 24575
$ws.Range("D108").Value = 5032
$ws.Range("E108").Value = 45379780
$ws.Range("C150").Value = 62801
$ws.Range("E150").Value = 119773069
$ws.Range("C151").Value = 23219
$ws.Range("E151").Value = 41385836
$ws.Range("C152").Value = 59474
$ws.Range("D152").Value = 12686
$ws.Range("E152").Value = 194367147
$ws.Range("C156").Value = 29368
$ws.Range("D156").Value = 9441
$ws.Range("E156").Value = 54626122
$ws.Range("C160").Value = 24092
$ws.Range("D160").Value = 4911
$ws.Range("E160").Value = 46768542
$ws.Range("C161").Value = 46661
$ws.Range("D161").Value = 13740
$ws.Range("E161").Value = 75831908
$ws.Range("C165").Value = 1283
$ws.Range("E165").Value = 2305108
$ws.Range("C167").Value = 180682
$ws.Range("D167").Value = 54468
$ws.Range("E167").Value = 392766708
$ws.Range("C168").Value = 351132
$ws.Range("D168").Value = 58991
$ws.Range("E168").Value = 691770856
$ws.Range("C169").Value = 174356
$ws.Range("D169").Value = 35290
$ws.Range("E169").Value = 675174361
$ws.Range("C170").Value = 67083
$ws.Range("D170").Value = 18257
$ws.Range("E170").Value = 135501701
$ws.Range("C172").Value = 33187
$ws.Range("D172").Value = 10605
$ws.Range("E172").Value = 66876948
$ws.Range("C173").Value = 219777
$ws.Range("D173").Value = 64858
$ws.Range("E173").Value = 391363321
$ws.Range("C174").Value = 70714
$ws.Range("D174").Value = 16552
$ws.Range("E174").Value = 160074078
$ws.Range("C176").Value = 62519
$ws.Range("E176").Value = 94836715
$ws.Range("D177").Value = 28523
$ws.Range("E177").Value = 89793823
$ws.Range("C178").Value = 132470
$ws.Range("D178").Value = 26897
$ws.Range("E178").Value = 254512122
$ws.Range("C179").Value = 97205
$ws.Range("E179").Value = 187019902
$ws.Range("D236").Value = 12898
$ws.Range("E236").Value = 79128086
$ws.Range("C238").Value = 40166
$ws.Range("D238").Value = 8874
$ws.Range("E238").Value = 132362267
$ws.Range("C242").Value = 17251
$ws.Range("D242").Value = 5606
$ws.Range("E242").Value = 34062405
$ws.Range("C243").Value = 8151
$ws.Range("D243").Value = 2479
$ws.Range("E243").Value = 15010881
$ws.Range("C246").Value = 15579
$ws.Range("E246").Value = 29806442
$ws.Range("C250").Value = 23934
$ws.Range("E250").Value = 34738409
$ws.Range("C259").Value = 15844
$ws.Range("E259").Value = 26808381
$ws.Range("C263").Value = 29248
$ws.Range("E263").Value = 42698696
$ws.Range("C264").Value = 38371
$ws.Range("E264").Value = 50552851
$ws.Range("C265").Value = 38813
$ws.Range("E265").Value = 70668096
$ws.Range("C266").Value = 65057
$ws.Range("D266").Value = 19560
$ws.Range("E266").Value = 107546151
$ws.Range("C267").Value = 18566
$ws.Range("D267").Value = 6932
$ws.Range("E267").Value = 30257719
$ws.Range("E273").Value = 215360795
$ws.Range("C274").Value = 24289
$ws.Range("E274").Value = 44907030
$ws.Range("C278").Value = 18994
$ws.Range("D278").Value = 6211
$ws.Range("E278").Value = 31951810
$ws.Range("C279").Value = 58092
$ws.Range("D279").Value = 18617
$ws.Range("E279").Value = 110316132
$ws.Range("C280").Value = 27515
$ws.Range("D280").Value = 8542
$ws.Range("E280").Value = 48878014
$ws.Range("C282").Value = 36569
$ws.Range("E282").Value = 51628872
$ws.Range("E283").Value = 56264203
$ws.Range("C285").Value = 67491
$ws.Range("D285").Value = 19748
$ws.Range("E285").Value = 105938913
$ws.Range("C298").Value = 26330
$ws.Range("E298").Value = 52089042
$ws.Range("C299").Value = 9665
$ws.Range("D299").Value = 3195
$ws.Range("E299").Value = 19424628
$ws.Range("C309").Value = 63357
$ws.Range("D309").Value = 25303
$ws.Range("E309").Value = 96928787
$ws.Range("C311").Value = 45151
$ws.Range("D311").Value = 8884
$ws.Range("E311").Value = 91548855
$ws.Range("D312").Value = 22143
$ws.Range("C313").Value = 12760
$ws.Range("E313").Value = 22268703
$ws.Range("C316").Value = 63111
$ws.Range("D316").Value = 19882
$ws.Range("E316").Value = 117193372
$ws.Range("C317").Value = 35685
$ws.Range("D317").Value = 10326
$ws.Range("E317").Value = 69229417
$ws.Range("C319").Value = 39654
$ws.Range("D319").Value = 9825
$ws.Range("E319").Value = 58440366
$ws.Range("C320").Value = 39310
$ws.Range("E320").Value = 52180782
$ws.Range("C321").Value = 44102
$ws.Range("D321").Value = 8992
$ws.Range("E321").Value = 86109259
$ws.Range("C322").Value = 67093
$ws.Range("E322").Value = 110653032
